# Auto-generated edit script applying numeric corrections described in the diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 4354585
$ws.Range("J17").Value = 4354585
$ws.Range("L17").Value = 13063755
$ws.Range("N17").Value = -13064091

# Row 121
$ws.Range("H121").Value = 4908.25
$ws.Range("I121").Value = 550
$ws.Range("J121").Value = 5243.5
$ws.Range("K121").Value = 1650
$ws.Range("L121").Value = 15730.5
$ws.Range("M121").Value = 97
$ws.Range("N121").Value = -19224.5

# Row 137
$ws.Range("H137").Value = 116786.34
$ws.Range("I137").Value = 135421.4
$ws.Range("J137").Value = 4976
$ws.Range("K137").Value = 406264.2
$ws.Range("L137").Value = 14928
$ws.Range("M137").Value = -403714.2
$ws.Range("N137").Value = -20028

# Row 138
$ws.Range("H138").Value = 4380.9697
$ws.Range("I138").Value = 2499.75
$ws.Range("J138").Value = 4640.448
$ws.Range("K138").Value = 7499.25
$ws.Range("L138").Value = 13921.344
$ws.Range("M138").Value = -2359.25
$ws.Range("N138").Value = -24201.344

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 14562.677
$ws.Range("I32").Value = 10901.148
$ws.Range("J32").Value = 28685.715
$ws.Range("K32").Value = 10901.148
$ws.Range("L32").Value = 28685.715
$ws.Range("M32").Value = -10614.148
$ws.Range("N32").Value = -29259.715

# Row 45
$ws.Range("H45").Value = 4289
$ws.Range("I45").Value = 4644.4443
$ws.Range("J45").Value = 3933.5557
$ws.Range("K45").Value = 4644.4443
$ws.Range("L45").Value = 3933.5557
$ws.Range("M45").Value = -4267.4443
$ws.Range("N45").Value = -4687.5557

# Row 102
$ws.Range("H102").Value = 986.3333
$ws.Range("I102").Value = 939.63635
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 939.63635
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 682.36365
$ws.Range("N102").Value = -4744

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 854.4722
$ws.Range("I94").Value = 377.81482
$ws.Range("J94").Value = 2284.4443
$ws.Range("K94").Value = 377.81482
$ws.Range("L94").Value = 2284.4443
$ws.Range("M94").Value = 73.18518
$ws.Range("N94").Value = -3186.4443

# Row 99
$ws.Range("H99").Value = 1528.1
$ws.Range("I99").Value = 1067.5
$ws.Range("J99").Value = 1835.1666
$ws.Range("K99").Value = 1067.5
$ws.Range("L99").Value = 1835.1666
$ws.Range("M99").Value = 430.5
$ws.Range("N99").Value = -4831.1666

# Row 113
$ws.Range("H113").Value = 2933.3333
$ws.Range("I113").Value = 2933.3333
$ws.Range("K113").Value = 2933.3333
$ws.Range("M113").Value = -763.3332999999998

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 53.875
$ws.Range("I7").Value = 47.75
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 47.75
$ws.Range("L7").Value = 60
$ws.Range("M7").Value = 65.25
$ws.Range("N7").Value = -286

# Row 31
$ws.Range("H31").Value = 5991.706
$ws.Range("I31").Value = 3853.5789
$ws.Range("J31").Value = 8700
$ws.Range("K31").Value = 3853.5789
$ws.Range("L31").Value = 8700
$ws.Range("M31").Value = -3558.5789
$ws.Range("N31").Value = -9290

# Row 34
$ws.Range("H34").Value = 5991.706
$ws.Range("I34").Value = 3853.5789
$ws.Range("J34").Value = 8700
$ws.Range("K34").Value = 3853.5789
$ws.Range("L34").Value = 8700
$ws.Range("M34").Value = -3651.5789
$ws.Range("N34").Value = -9104

$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 71.666664
$ws.Range("I38").Value = 57.8
$ws.Range("K38").Value = 173.4
$ws.Range("M38").Value = 173.6

# Row 107
$ws.Range("H107").Value = 6287.2812
$ws.Range("I107").Value = 8184.125
$ws.Range("J107").Value = 596.75
$ws.Range("K107").Value = 24552.375
$ws.Range("L107").Value = 1790.25
$ws.Range("M107").Value = -22632.375
$ws.Range("N107").Value = -5630.25

# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# Row 113
$ws.Range("H113").Value = 807.8333
$ws.Range("I113").Value = 615.8333
$ws.Range("J113").Value = 871.8333
$ws.Range("K113").Value = 1847.4999
$ws.Range("L113").Value = 2615.4999
$ws.Range("M113").Value = 322.5001
$ws.Range("N113").Value = -6955.4999

# Row 114
$ws.Range("H114").Value = 3270
$ws.Range("I114").Value = 3360
$ws.Range("J114").Value = 3000
$ws.Range("K114").Value = 10080
$ws.Range("L114").Value = 9000
$ws.Range("M114").Value = -6826
$ws.Range("N114").Value = -15508

# Row 116
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

# Row 131
$ws.Range("H131").Value = 753.77
$ws.Range("J131").Value = 780.20654
$ws.Range("L131").Value = 2340.61962
$ws.Range("N131").Value = -12420.61962

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 15333.333
$ws.Range("J43").Value = 26000
$ws.Range("L43").Value = 26000
$ws.Range("N43").Value = -26302

# Row 46
$ws.Range("H46").Value = 39999
$ws.Range("J46").Value = 39999
$ws.Range("L46").Value = 39999
$ws.Range("N46").Value = -40311

# Row 70
$ws.Range("H70").Value = 3129080.5
$ws.Range("I70").Value = 4091.3333
$ws.Range("K70").Value = 4091.3333
$ws.Range("M70").Value = -3821.3333

# Row 73
$ws.Range("H73").Value = 3129080.5
$ws.Range("I73").Value = 4091.3333
$ws.Range("K73").Value = 4091.3333
$ws.Range("M73").Value = -3155.3333

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 1996.5714
$ws.Range("I82").Value = 1300
$ws.Range("J82").Value = 2925.3333
$ws.Range("K82").Value = 1300
$ws.Range("L82").Value = 2925.3333
$ws.Range("M82").Value = -939
$ws.Range("N82").Value = -3647.3333

# Row 85
$ws.Range("H85").Value = 1996.5714
$ws.Range("I85").Value = 1300
$ws.Range("J85").Value = 2925.3333
$ws.Range("K85").Value = 1300
$ws.Range("L85").Value = 2925.3333
$ws.Range("M85").Value = -52
$ws.Range("N85").Value = -5421.3333

# Row 93
$ws.Range("H93").Value = 965.1795
$ws.Range("I93").Value = 896.9143
$ws.Range("J93").Value = 1562.5
$ws.Range("K93").Value = 896.9143
$ws.Range("L93").Value = 1562.5
$ws.Range("M93").Value = 351.0857
$ws.Range("N93").Value = -4058.5

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 700
$ws.Range("J100").Value = 700
$ws.Range("L100").Value = 1400
$ws.Range("N100").Value = -2482
